$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.5
$ws.Range("N2").Value = 7.5
$ws.Range("U2").Value = 2.05
$ws.Range("V2").Value = 1.7
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 7.5
$ws.Range("AC2").Value = 7.5
$ws.Range("AD2").Value = 6.5
$ws.Range("AE2").Value = 19
$ws.Range("AN2").Value = 3.6
